# Added Jimmy Jordan Abstract
#
# Row 24/25 on Sheet1 both describe event #12 (October 4th), which has two
# speakers. Jimmy Jordan's talk now gets its real abstract (replacing the
# "Details coming soon!" placeholder), and his row is moved to be first
# (row 24), pushing Adam Howell's row down to row 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$jimmyAbstract = "The newspaper industry in the U.S. has been undergoing a significant decline since the early 2000s. The rise of the internet ushered in an era of information abundance and accessibility. However, this transition has severely impacted the traditional business models that once supported local newspapers, resulting in noticeable gaps in community coverage across the country. All of this is set against a climate where trust in media is approaching historic lows in the U.S.<br><br>This presentation will illustrate these challenges, both on a local scale and within a broader context. Additionally, we will delve into the new models that are being explored in an effort to revitalize local reporting. "

# Capture the current (pre-edit) contents of both speaker rows.
$adamSpeaker  = $ws.Range("D24").Value()
$adamPosition = $ws.Range("E24").Value()
$adamFlavor   = $ws.Range("F24").Value()
$adamTitle    = $ws.Range("G24").Value()
$adamAbstract = $ws.Range("H24").Value()

$jimmySpeaker  = $ws.Range("D25").Value()
$jimmyPosition = $ws.Range("E25").Value()
$jimmyFlavor   = $ws.Range("F25").Value()
$jimmyTitle    = $ws.Range("G25").Value()

# Row 24 becomes Jimmy Jordan, with his new full abstract.
$ws.Range("D24").Value = $jimmySpeaker
$ws.Range("E24").Value = $jimmyPosition
$ws.Range("F24").Value = $jimmyFlavor
$ws.Range("G24").Value = $jimmyTitle
$ws.Range("H24").Value = $jimmyAbstract

# Row 25 becomes Adam Howell (unchanged content, just moved down a row).
$ws.Range("D25").Value = $adamSpeaker
$ws.Range("E25").Value = $adamPosition
$ws.Range("F25").Value = $adamFlavor
$ws.Range("G25").Value = $adamTitle
$ws.Range("H25").Value = $adamAbstract

# Match the saved sheet view/selection state from the edit.
$ws.Activate()
$ws.Range("F28").Select()
$excel.ActiveWindow.ScrollRow = 13
